$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the now-obsolete "60th"/"70th" Ground Truth Point notes;
# these ground truth milestones are being superseded.
$ws.Range("E68").ClearContents()
$ws.Range("E78").ClearContents()

# Row 97 gets a new note instead of a milestone marker.
$ws.Range("E97").Value = "Way too close to top"

# The "90th Ground Truth Point" marker moves from row 98 to row 101.
$ws.Range("E98").ClearContents()
$ws.Range("E101").Value = "90th Ground Truth Point"

# Fill in newly measured points for rows 91-96, 98-100, 102 (B:D),
# and 101 (B:D) as part of continuing to 91 ground truth points.
$ws.Range("B91").Value = 32
$ws.Range("C91").Value = 56
$ws.Range("D91").Value = 25

$ws.Range("B92").Value = 41
$ws.Range("C92").Value = 50
$ws.Range("D92").Value = 27

$ws.Range("B93").Value = 40
$ws.Range("C93").Value = 47
$ws.Range("D93").Value = 24

$ws.Range("B94").Value = 30
$ws.Range("C94").Value = 49
$ws.Range("D94").Value = 22

$ws.Range("B95").Value = 26
$ws.Range("C95").Value = 46
$ws.Range("D95").Value = 20

$ws.Range("B96").Value = 30
$ws.Range("C96").Value = 26
$ws.Range("D96").Value = 20

$ws.Range("B98").Value = 26
$ws.Range("C98").Value = 34
$ws.Range("D98").Value = 20

$ws.Range("B99").Value = 33
$ws.Range("C99").Value = 55
$ws.Range("D99").Value = 21

$ws.Range("B100").Value = 31
$ws.Range("C100").Value = 61
$ws.Range("D100").Value = 25

$ws.Range("B101").Value = 20
$ws.Range("C101").Value = 55
$ws.Range("D101").Value = 15

$ws.Range("B102").Value = 37
$ws.Range("C102").Value = 57
$ws.Range("D102").Value = 26

# Update the visible top-left cell / selection to reflect scrolling to the
# newly added rows.
$ws.Range("B103").Select()
$excel.ActiveWindow.ScrollRow = 94
